# Add a new "Linked_Poster" column and a new job posting row (JD_003)
# per commit "Add Job Posting with Job_Id=JD_003".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: I1 = Linked_Poster, formatted like the other headers ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(1, 9).Value = "Linked_Poster"

# --- New row 4: JD_003 / Senior Dotnet Engineer -------------------------
$ws.Cells.Item(4, 1).Value = "JD_003"
$ws.Cells.Item(4, 2).Value = "Senior Dotnet Engineer"
$ws.Cells.Item(4, 3).Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0

# Re-fit the new row's height (the multi-line description otherwise leaves a
# stale custom row height behind).
$ws.Rows.Item(4).AutoFit()
